$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.675.92"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.070.72"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'232.77"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'58.40"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("D12").Value = "2.375.70"
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").Value = "'20.90"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'0.772"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "'5.36"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "2.078.81"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "37.634.58"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "'71.20"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'227.97"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "'171.15"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "'0.136"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").Value = "'19.44"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "'4.67"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "'4.64"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "'99.61"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "'0.0971"
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +6.32%  "
$ws.Range("D45").Value = "1.434.33"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "'4.17"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "'7.38"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "2.260.56"
$ws.Range("E51").Value = "  -2.25%  "
